$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.760.30"
$ws.Range("E2").Value = "  -7.30%  "
$ws.Range("D3").Value = "2.536.59"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'296.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.24%  "
$ws.Range("D6").Value = "'93.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("D7").Value = "'0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.13%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.97%  "
$ws.Range("D10").Value = "'35.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.20%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("D12").Value = "'7.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.75%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.923.38"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.107"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "2.533.99"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "'0.862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.82%  "
$ws.Range("D17").Value = "'14.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.96%  "
$ws.Range("D18").Value = "42.797.97"
$ws.Range("E18").Value = "  -7.44%  "
$ws.Range("D19").Value = "0.0₃0962"
$ws.Range("E19").Value = "  -4.98%  "
$ws.Range("D20").Value = "'6.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").Value = "'12.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "'71.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'257.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.32%  "
$ws.Range("D24").Value = "'2.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("D26").Value = "'29.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.86%  "
$ws.Range("D30").Value = "'35.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.84%  "
$ws.Range("E31").Value = "  -6.27%  "
$ws.Range("D32").Value = "'150.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").Value = "'3.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.84%  "
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("D35").Value = "'2.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").Value = "'0.0792"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("D37").Value = "'0.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.13%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "'24.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.61%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("D40").Value = "'15.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").Value = "'3.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").Value = "'0.0307"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.08%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.066.13"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'3.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'84.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.77%  "
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").Value = "'8.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.21%  "
$ws.Range("D49").Value = "2.781.14"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").Value = "'1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("D51").Value = "'102.67"
$ws.Range("D51").Style = "Normal"
